# Commit: "driver changed to chrome"
# Populates the itemTotal/Tax/Total columns (G:I) for the two order rows
# (rows 2 and 3) with their dollar-formatted text values, e.g. "$55.97".
#
# The values must land as plain text (shared-string) cells, matching the
# source data which stores currency amounts as strings like "$55.97"
# rather than numeric currency cells. Assigning a "$"-prefixed string
# straight to .Value makes Excel auto-coerce it into a numeric currency
# cell, so each cell is briefly switched to Text number format while the
# value is entered, then the formatting is cleared again so the cell is
# left with plain/default formatting and a text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "G2"; Value = "$55.97" },
    @{ Cell = "H2"; Value = "$4.48" },
    @{ Cell = "I2"; Value = "$60.45" },
    @{ Cell = "G3"; Value = "$25.98" },
    @{ Cell = "H3"; Value = "$2.08" },
    @{ Cell = "I3"; Value = "$28.06" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.ClearFormats()
}
